$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.071.72'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '1.623.34'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.00'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.515'
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0628'
$ws.Range("E8").Value = '  +0.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.250'
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.89'
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0845'
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").Value = '1.852.20'
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("D13").Value = '1.612.88'
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.13'
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.539'
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '27.062.51'
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.34'
$ws.Range("E17").Value = '  -3.38%  '
$ws.Range("D18").Value = '0.0₃0735'
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.14'
$ws.Range("E19").Value = '  -1.17%  '
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.83'
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.34'
$ws.Range("E22").Value = '  -1.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.33'
$ws.Range("E23").Value = '  -7.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.02'
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.15'
$ws.Range("E25").Value = '  +1.28%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.35'
$ws.Range("E27").Value = '  -0.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.115'
$ws.Range("E28").Value = '  -3.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.51'
$ws.Range("E29").Value = '  -0.93%  '
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("E31").Value = '  -1.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.707'
$ws.Range("E33").Value = '  +30.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.99'
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("D35").Value = '1.354.95'
$ws.Range("E35").Value = '  +3.61%  '
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0176'
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.842'
$ws.Range("E39").Value = '  -1.66%  '
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.23'
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.798'
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.35'
$ws.Range("E43").Value = '  +1.24%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.93'
$ws.Range("E44").Value = '  +3.83%  '
$ws.Range("D45").Value = '1.763.19'
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.65'
$ws.Range("E46").Value = '  +3.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.84'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.861'
$ws.Range("E48").Value = '  +28.67%  '
$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("E49").Value = '  -1.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.100'
$ws.Range("E50").Value = '  +4.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0512'
$ws.Range("E51").Value = '  +0.13%  '

Write-Host "Updated cryptos list"
